{"js": "// Update the worksheet date and the 25 division problems in the table,\n// matching the target edit exactly (old -> new, in document order).\nconst replacements = [\n  [\"2023-07-26 Wednesday\", \"2023-07-27 Thursday\"],\n  [\"39\u00f78=\", \"29\u00f77=\"],\n  [\"57\u00f74=\", \"75\u00f78=\"],\n  [\"56\u00f77=\", \"57\u00f77=\"],\n  [\"28\u00f75=\", \"88\u00f76=\"],\n  [\"57\u00f74=\", \"45\u00f78=\"],\n  [\"91\u00f73=\", \"70\u00f72=\"],\n  [\"66\u00f79=\", \"73\u00f76=\"],\n  [\"35\u00f73=\", \"26\u00f78=\"],\n  [\"26\u00f73=\", \"20\u00f77=\"],\n  [\"40\u00f73=\", \"72\u00f72=\"],\n  [\"28\u00f73=\", \"52\u00f75=\"],\n  [\"42\u00f77=\", \"92\u00f79=\"],\n  [\"84\u00f77=\", \"88\u00f74=\"],\n  [\"80\u00f72=\", \"61\u00f79=\"],\n  [\"30\u00f75=\", \"29\u00f77=\"],\n  [\"19\u00f75=\", \"33\u00f75=\"],\n  [\"95\u00f73=\", \"88\u00f74=\"],\n  [\"30\u00f76=\", \"37\u00f78=\"],\n  [\"22\u00f73=\", \"51\u00f77=\"],\n  [\"80\u00f73=\", \"46\u00f73=\"],\n  [\"25\u00f76=\", \"61\u00f74=\"],\n  [\"58\u00f78=\", \"31\u00f73=\"],\n  [\"71\u00f72=\", \"38\u00f76=\"],\n  [\"77\u00f75=\", \"96\u00f75=\"],\n  [\"14\u00f72=\", \"76\u00f77=\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet r = 0;\nfor (let i = 0; i < paragraphs.items.length && r < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  const [oldText, newText] = replacements[r];\n  if (para.text === oldText) {\n    para.insertText(newText, Word.InsertLocation.replace);\n    r++;\n  }\n}\nawait context.sync();\n\nif (r !== replacements.length) {\n  throw new Error(\n    `Only applied ${r} of ${replacements.length} replacements`\n  );\n}\n", "ps1": "# Update the worksheet date and the 25 division problems in the table,\n# matching the target edit exactly (old -> new, by table position).\n$d = $word.ActiveDocument\n\n# Word's COM Range.Text includes trailing control characters (paragraph\n# mark \\r for normal paragraphs, \\r\\a cell-end mark for table cells), so\n# strip those before comparing against the plain text we expect.\nfunction TrimMarks([string]$s) {\n    return $s.TrimEnd([char]13, [char]7)\n}\n\n# Heading paragraph with the date.\n$titleOld = \"2023-07-26 Wednesday\"\n$titleNew = \"2023-07-27 Thursday\"\n$p1 = $d.Paragraphs(1)\nif ((TrimMarks $p1.Range.Text) -eq $titleOld) {\n    $p1.Range.Text = $titleNew\n}\n\n# Division problems live in table 1; data rows are 1,5,9,13,17 (1-based),\n# each with 5 populated columns (25 cells total). The old/new pairs below\n# are listed in (row,col) document order, matching the source diff exactly\n# -- note some old values repeat with different replacements (e.g. \"57\u00f74=\"\n# becomes \"75\u00f78=\" in one cell and \"45\u00f78=\" in another), so positional\n# (row,col) targeting is used instead of a global text search/replace.\n$table = $d.Tables(1)\n$dataRows = @(1, 5, 9, 13, 17)\n\n$oldValues = @(\n    \"39\u00f78=\", \"57\u00f74=\", \"56\u00f77=\", \"28\u00f75=\", \"57\u00f74=\",\n    \"91\u00f73=\", \"66\u00f79=\", \"35\u00f73=\", \"26\u00f73=\", \"40\u00f73=\",\n    \"28\u00f73=\", \"42\u00f77=\", \"84\u00f77=\", \"80\u00f72=\", \"30\u00f75=\",\n    \"19\u00f75=\", \"95\u00f73=\", \"30\u00f76=\", \"22\u00f73=\", \"80\u00f73=\",\n    \"25\u00f76=\", \"58\u00f78=\", \"71\u00f72=\", \"77\u00f75=\", \"14\u00f72=\"\n)\n$newValues = @(\n    \"29\u00f77=\", \"75\u00f78=\", \"57\u00f77=\", \"88\u00f76=\", \"45\u00f78=\",\n    \"70\u00f72=\", \"73\u00f76=\", \"26\u00f78=\", \"20\u00f77=\", \"72\u00f72=\",\n    \"52\u00f75=\", \"92\u00f79=\", \"88\u00f74=\", \"61\u00f79=\", \"29\u00f77=\",\n    \"33\u00f75=\", \"88\u00f74=\", \"37\u00f78=\", \"51\u00f77=\", \"46\u00f73=\",\n    \"61\u00f74=\", \"31\u00f73=\", \"38\u00f76=\", \"96\u00f75=\", \"76\u00f77=\"\n)\n\n$idx = 0\nforeach ($row in $dataRows) {\n    for ($col = 1; $col -le 5; $col++) {\n        $oldText = $oldValues[$idx]\n        $newText = $newValues[$idx]\n        $cell = $table.Cell($row, $col)\n        if ((TrimMarks $cell.Range.Text) -eq $oldText) {\n            $cell.Range.Text = $newText\n        }\n        $idx++\n    }\n}\n"}
